$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the typo "oof" -> "of" in "...during the development oof cancel
#    cells." ("oof" occurs exactly once in the whole document, so matching
#    the full phrase is unambiguous).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "development oof cancel cells.", $true, $false, $false, $false, $false,
    $true, 1, $false, "development of cancel cells.", 2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert three new paragraphs ("Synthetic biology", "DNA libraries", and
#    "Nicking enzymes") right after the "Markers of epigenetic modifications"
#    paragraph we just fixed, and before the blank paragraph that follows it.
# ---------------------------------------------------------------------------

# Re-locate the (now-fixed) paragraph by its text so the insertion point is
# found robustly rather than via a hard-coded paragraph index.
$findRange = $d.Content.Duplicate
$findRange.Find.Execute(
    "development of cancel cells.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null

$targetParagraph = $findRange.Paragraphs(1)
$targetRange = $targetParagraph.Range
$targetRange.Collapse(0)

# Create a fresh, empty paragraph right after the target paragraph. We then
# fill *that* empty paragraph's range (paragraph mark included) with our
# fully-formatted OOXML, since InsertXML replaces the contents of whatever
# range it is invoked on (calling it directly on the original paragraph's
# collapsed end would clobber that paragraph instead of inserting after it).
$targetRange.InsertParagraphAfter()
$newRange = $targetParagraph.Next().Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Synthetic biology</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>: aims to redesign</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">, or create </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>biological systems</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">to enhance existing </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>or create</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> new abilities. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">Some DNA assembly technologies are based on </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>R</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>E</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>ase</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>-ligation methodology</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">the ability of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>REases</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> to cleave DNA outside of the recognition sequence</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> or exonucleases activity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>DNA libraries: the use of type IIS enzymes as tagging enzymes have allow the identification and quantification of a large numbers o mRNA transcripts</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">. In addition, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>REases</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Acul</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> and </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>USER</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">) have </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>been used to insert tags into sample DNA, which is the amplified to form long, single-stranded DNA &#8220;nanoballs&#8221; that are templates in chip-based sequencing-by-litigation methodology.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Nick</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>ing enzymes</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">: Some type IIS enzymes have two different catalytic sites. By inactivating one catalytic site at a time, these enzymes are converted strand-specific DNA nicking enzymes (cleaving only one DNA strand) and are useful molecular investigating tools for </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>research</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> on altering </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>DNA, generation</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> oof reporter plasmids with modified structures, and the creation of a DNA motor that transports a DNA payload.</w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xml)

Write-Output "done"
